$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the cryptos list (price + 1h volume change columns) with the latest
# scrape. Price cells (column D) are stored as text in the sheet, so
# numeric-looking values are prefixed with a leading apostrophe to force
# Excel to keep them as text (preserving formatting like trailing zeros)
# instead of auto-converting them to numbers.
$ws.Range('D2').Value = '26.357.61'
$ws.Range('E2').Value = '  +0.50%  '
$ws.Range('D3').Value = '1.592.22'
$ws.Range('E3').Value = '  +0.53%  '
$ws.Range('E4').Value = '  -0.29%  '
$ws.Range('D5').Value = '''211.67'
$ws.Range('E5').Value = '  +0.87%  '
$ws.Range('D6').Value = '''0.504'
$ws.Range('E6').Value = '  -0.15%  '
$ws.Range('E8').Value = '  +0.16%  '
$ws.Range('E9').Value = '  -0.12%  '
$ws.Range('D10').Value = '''19.46'
$ws.Range('E10').Value = '  -0.53%  '
$ws.Range('D11').Value = '''0.0847'
$ws.Range('E11').Value = '  +0.14%  '
$ws.Range('D12').Value = '1.815.66'
$ws.Range('E12').Value = '  +0.55%  '
$ws.Range('D13').Value = '1.627.22'
$ws.Range('E13').Value = '  +3.48%  '
$ws.Range('D14').Value = '''4.05'
$ws.Range('E14').Value = '  +0.74%  '
$ws.Range('E15').Value = '  +0.93%  '
$ws.Range('D16').Value = '''64.64'
$ws.Range('E16').Value = '  -0.06%  '
$ws.Range('D17').Value = '26.355.88'
$ws.Range('E17').Value = '  +0.50%  '
$ws.Range('D18').Value = '0.0₃0733'
$ws.Range('E18').Value = '  -0.78%  '
$ws.Range('D19').Value = '''7.50'
$ws.Range('E19').Value = '  +3.72%  '
$ws.Range('D20').Value = '''212.50'
$ws.Range('E20').Value = '  +2.65%  '
$ws.Range('E21').Value = '  -0.29%  '
$ws.Range('E22').Value = '  +0.99%  '
$ws.Range('E23').Value = '  -1.12%  '
$ws.Range('E24').Value = '  +1.90%  '
$ws.Range('D25').Value = '''144.91'
$ws.Range('E25').Value = '  +0.12%  '
$ws.Range('E26').Value = '  -0.29%  '
$ws.Range('E27').Value = '  +0.43%  '
$ws.Range('E28').Value = '  -0.56%  '
$ws.Range('E29').Value = '  -0.27%  '
$ws.Range('E30').Value = '  +0.03%  '
$ws.Range('E31').Value = '  +0.94%  '
$ws.Range('E32').Value = '  -0.26%  '
$ws.Range('D33').Value = '''2.99'
$ws.Range('E33').Value = '  +1.35%  '
$ws.Range('D34').Value = '1.344.05'
$ws.Range('E34').Value = '  +4.21%  '
$ws.Range('E35').Value = '  -1.03%  '
$ws.Range('D36').Value = '''0.604'
$ws.Range('E36').Value = '  -0.53%  '
$ws.Range('D37').Value = '''1.49'
$ws.Range('E37').Value = '  +0.12%  '
$ws.Range('E38').Value = '  -0.06%  '
$ws.Range('E39').Value = '  -17.46%  '
$ws.Range('D40').Value = '''0.820'
$ws.Range('E40').Value = '  +0.47%  '
$ws.Range('D41').Value = '''5.80'
$ws.Range('E41').Value = '  +4.41%  '
$ws.Range('E42').Value = '  -0.26%  '
$ws.Range('E43').Value = '  +0.24%  '
$ws.Range('E44').Value = '  -0.68%  '
$ws.Range('D45').Value = '1.728.30'
$ws.Range('E45').Value = '  +0.53%  '
$ws.Range('D46').Value = '''61.72'
$ws.Range('E46').Value = '  -1.04%  '
$ws.Range('D47').Value = '''87.92'
$ws.Range('E47').Value = '  -0.97%  '
$ws.Range('D48').Value = '0.0₆0105'
$ws.Range('E48').Value = '  +6.65%  '
$ws.Range('E49').Value = '  -2.92%  '
$ws.Range('D50').Value = '''0.0990'
$ws.Range('E50').Value = '  -2.89%  '
$ws.Range('E51').Value = '  -0.77%  '
